$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G task notes added to existing rows (29-39) ---
$ws.Range("G29").Value = "performance tweeks, view subject , subject attendance"
$ws.Range("G30").Value = "minimizer and profile pic upload"
$ws.Range("G31").Value = "GUI, profile pic for students"
$ws.Range("G32").Value = "preferences for combined timetable, disable on labs completed"
$ws.Range("G33").Value = "bug fixes, add student, unselectable text, "
$ws.Range("G34").Value = "fixed new week redundancy"
$ws.Range("G36").Value = "time based timetable edit, fix date printing"
$ws.Range("G37").Value = "convert full site to ajax"
$ws.Range("G38").Value = "convert full site to ajax"
$ws.Range("G39").Value = "convert full site to ajax"

# --- Updated hour values ---
$ws.Range("B34").Value = 5
$ws.Range("B35").Value = 0
$ws.Range("B37").Value = 12
$ws.Range("B38").Value = 5

# --- Column G notes for rows 40-42 (existing rows, new cells) ---
$ws.Range("G40").Value = "convert full site to ajax"
$ws.Range("G41").Value = "convert full site to ajax"
$ws.Range("G42").Value = "convert full site to ajax"

# --- New row 44 spacer with styled (but empty) G44 cell ---
$ws.Range("G43").Style = $ws.Range("G39").Style
$ws.Range("G44").Style = $ws.Range("G39").Style

# --- Row 45: weekly summary continuation ---
$ws.Range("A45").Value = 4
$ws.Range("B45").Value = 6
$ws.Range("G45").Value = "convert full site to ajax"
$ws.Range("G45").Style = $ws.Range("G39").Style

# --- New rows 46-53: weekly summary + task notes ---
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = 5
$ws.Range("G46").Value = "ajax bug fix"

$ws.Range("A47").Value = 6
$ws.Range("B47").Value = 9
$ws.Range("G47").Value = "subject selection bugs, slot management"

$ws.Range("A48").Value = 7
$ws.Range("B48").Value = 9
$ws.Range("G48").Value = "new gui, session management"

$ws.Range("A49").Value = 8
$ws.Range("B49").Value = 7
$ws.Range("G49").Value = "gui bugs, data tables"

$ws.Range("A50").Value = 9
$ws.Range("B50").Value = 9
$ws.Range("G50").Value = "save stud details"

$ws.Range("A51").Value = 10
$ws.Range("B51").Value = 5
$ws.Range("G51").Value = "save stud details , back button "

$ws.Range("A52").Value = 11
$ws.Range("B52").Value = 8
$ws.Range("G52").Value = "subject batch for students, back button bug fix"

$ws.Range("A53").Value = 12
$ws.Range("B53").Value = 5
$ws.Range("G53").Value = "soldering rpi keypad and fit inside box"

$ws.Range("A54").Value = 13

# Apply matching styles for new A/B cells (same look as rows 40-42)
$ws.Range("A45:A54").Style = $ws.Range("A42").Style
$ws.Range("B45:B54").Style = $ws.Range("B42").Style
$ws.Range("G46:G53").Style = $ws.Range("B42").Style

# --- View / pane changes ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A29").Select()
$win.FreezePanes = $true
$ws.Range("H48").Select()
